$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.028588742700773
$ws.Cells.Item(2, 4).Value = 1.040369101480999
$ws.Cells.Item(2, 5).Value = 1.0497477325152
$ws.Cells.Item(2, 6).Value = 1.05431897958277
$ws.Cells.Item(2, 9).Value = 1.037405673360392
$ws.Cells.Item(2, 10).Value = 1.033740225454385
$ws.Cells.Item(2, 11).Value = 1.043151645741665
$ws.Cells.Item(2, 12).Value = 1.052503929197246
$ws.Cells.Item(2, 13).Value = 1.057062518321465
$ws.Cells.Item(2, 14).Value = 1.015265176555708

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.02943719451359
$ws.Cells.Item(3, 4).Value = 1.041042280958627
$ws.Cells.Item(3, 5).Value = 1.050611502892757
$ws.Cells.Item(3, 6).Value = 1.05514947515687
$ws.Cells.Item(3, 9).Value = 1.037573137457105
$ws.Cells.Item(3, 10).Value = 1.034229867792003
$ws.Cells.Item(3, 11).Value = 1.043635595653682
$ws.Cells.Item(3, 12).Value = 1.053179874900278
$ws.Cells.Item(3, 13).Value = 1.057706189284808
$ws.Cells.Item(3, 14).Value = 1.015427485604379

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.029986595296295
$ws.Cells.Item(4, 4).Value = 1.041477900858751
$ws.Cells.Item(4, 5).Value = 1.051171463629327
$ws.Cells.Item(4, 6).Value = 1.055687436286715
$ws.Cells.Item(4, 9).Value = 1.037679865864187
$ws.Cells.Item(4, 10).Value = 1.034546448861465
$ws.Cells.Item(4, 11).Value = 1.043948090984848
$ws.Cells.Item(4, 12).Value = 1.05361765714817
$ws.Cells.Item(4, 13).Value = 1.058122612231571
$ws.Cells.Item(4, 14).Value = 1.015532407264997

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.030217656656382
$ws.Cells.Item(5, 4).Value = 1.041661040338195
$ws.Cells.Item(5, 5).Value = 1.051407119131885
$ws.Cells.Item(5, 6).Value = 1.055913730927158
$ws.Cells.Item(5, 9).Value = 1.03772434290609
$ws.Cells.Item(5, 10).Value = 1.034679478439721
$ws.Cells.Item(5, 11).Value = 1.044079306408948
$ws.Cells.Item(5, 12).Value = 1.053801795076136
$ws.Cells.Item(5, 13).Value = 1.058297656851525
$ws.Cells.Item(5, 14).Value = 1.015576491213385

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.030256458302288
$ws.Cells.Item(6, 4).Value = 1.041691790496758
$ws.Cells.Item(6, 5).Value = 1.051446701218169
$ws.Cells.Item(6, 6).Value = 1.055951734701026
$ws.Cells.Item(6, 9).Value = 1.037731787809563
$ws.Cells.Item(6, 10).Value = 1.034701811073503
$ws.Cells.Item(6, 11).Value = 1.044101328773494
$ws.Cells.Item(6, 12).Value = 1.053832718127972
$ws.Cells.Item(6, 13).Value = 1.058327046418802
$ws.Cells.Item(6, 14).Value = 1.01558389161703

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.029989682384845
$ws.Cells.Item(7, 4).Value = 1.041480347960528
$ws.Cells.Item(7, 5).Value = 1.05117461149747
$ws.Cells.Item(7, 6).Value = 1.055690459515598
$ws.Cells.Item(7, 9).Value = 1.037680461708327
$ws.Cells.Item(7, 10).Value = 1.034548226650793
$ws.Cells.Item(7, 11).Value = 1.043949844912812
$ws.Cells.Item(7, 12).Value = 1.053620117239686
$ws.Cells.Item(7, 13).Value = 1.058124951265307
$ws.Cells.Item(7, 14).Value = 1.015532996416103

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.028875398144239
$ws.Cells.Item(8, 4).Value = 1.040596598824161
$ws.Cells.Item(8, 5).Value = 1.050039430593835
$ws.Cells.Item(8, 6).Value = 1.054599529281014
$ws.Cells.Item(8, 9).Value = 1.037462606047963
$ws.Cells.Item(8, 10).Value = 1.033905753479453
$ws.Cells.Item(8, 11).Value = 1.043315333185388
$ws.Cells.Item(8, 12).Value = 1.052732284431967
$ws.Cells.Item(8, 13).Value = 1.057280064559352
$ws.Cells.Item(8, 14).Value = 1.015320050761297

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.026914982427612
$ws.Cells.Item(9, 4).Value = 1.039039610680459
$ws.Cells.Item(9, 5).Value = 1.04804716932541
$ws.Cells.Item(9, 6).Value = 1.052681656573718
$ws.Cells.Item(9, 9).Value = 1.037066254795009
$ws.Cells.Item(9, 10).Value = 1.032771769349525
$ws.Cells.Item(9, 11).Value = 1.042192304635835
$ws.Cells.Item(9, 12).Value = 1.051170942098223
$ws.Cells.Item(9, 13).Value = 1.05579075996786
$ws.Cells.Item(9, 14).Value = 1.014944042010978

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.025610206293133
$ws.Cells.Item(10, 4).Value = 1.038001918925971
$ws.Cells.Item(10, 5).Value = 1.046724530486361
$ws.Cells.Item(10, 6).Value = 1.051406201157186
$ws.Cells.Item(10, 9).Value = 1.036793692742658
$ws.Cells.Item(10, 10).Value = 1.032014596685644
$ws.Cells.Item(10, 11).Value = 1.041440380095267
$ws.Cells.Item(10, 12).Value = 1.050132244151938
$ws.Cells.Item(10, 13).Value = 1.054797643411545
$ws.Cells.Item(10, 14).Value = 1.014692875987841

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.02504575594942
$ws.Cells.Item(11, 4).Value = 1.037552680159951
$ws.Cells.Item(11, 5).Value = 1.046153148828079
$ws.Cells.Item(11, 6).Value = 1.050854679572125
$ws.Cells.Item(11, 9).Value = 1.036673705650743
$ws.Cells.Item(11, 10).Value = 1.031686466581431
$ws.Cells.Item(11, 11).Value = 1.041114037199408
$ws.Cells.Item(11, 12).Value = 1.049683016612916
$ws.Cells.Item(11, 13).Value = 1.054367574024509
$ws.Cells.Item(11, 14).Value = 1.014584006511253

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.024836174455366
$ws.Cells.Item(12, 4).Value = 1.037385827951234
$ws.Cells.Item(12, 5).Value = 1.045941113834292
$ws.Cells.Item(12, 6).Value = 1.05064993579003
$ws.Cells.Item(12, 9).Value = 1.036628842561584
$ws.Cells.Item(12, 10).Value = 1.031564545147692
$ws.Cells.Item(12, 11).Value = 1.040992706887443
$ws.Cells.Item(12, 12).Value = 1.049516235486751
$ws.Cells.Item(12, 13).Value = 1.054207822151726
$ws.Cells.Item(12, 14).Value = 1.014543550990046

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.024881126731874
$ws.Cells.Item(13, 4).Value = 1.037421617606213
$ws.Cells.Item(13, 5).Value = 1.045986586908371
$ws.Cells.Item(13, 6).Value = 1.050693848752234
$ws.Cells.Item(13, 9).Value = 1.036638479161997
$ws.Cells.Item(13, 10).Value = 1.031590699469002
$ws.Cells.Item(13, 11).Value = 1.041018737699325
$ws.Cells.Item(13, 12).Value = 1.049552006870039
$ws.Cells.Item(13, 13).Value = 1.054242089669854
$ws.Cells.Item(13, 14).Value = 1.01455222957893

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.025028430227438
$ws.Cells.Item(14, 4).Value = 1.037538887786729
$ws.Cells.Item(14, 5).Value = 1.046135617826332
$ws.Cells.Item(14, 6).Value = 1.050837753012364
$ws.Cells.Item(14, 9).Value = 1.036670003257623
$ws.Cells.Item(14, 10).Value = 1.031676389312738
$ws.Cells.Item(14, 11).Value = 1.041104010275608
$ws.Cells.Item(14, 12).Value = 1.049669228751735
$ws.Cells.Item(14, 13).Value = 1.054354368976903
$ws.Cells.Item(14, 14).Value = 1.014580662780465

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.025119199481642
$ws.Cells.Item(15, 4).Value = 1.037611143864949
$ws.Cells.Item(15, 5).Value = 1.046227467461332
$ws.Cells.Item(15, 6).Value = 1.050926432588613
$ws.Cells.Item(15, 9).Value = 1.036689387288253
$ws.Cells.Item(15, 10).Value = 1.031729180480294
$ws.Cells.Item(15, 11).Value = 1.041156534719193
$ws.Cells.Item(15, 12).Value = 1.049741463931124
$ws.Cells.Item(15, 13).Value = 1.05442354735391
$ws.Cells.Item(15, 14).Value = 1.014598179233722

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.025647678213917
$ws.Cells.Item(16, 4).Value = 1.038031735440512
$ws.Cells.Item(16, 5).Value = 1.046762479375989
$ws.Cells.Item(16, 6).Value = 1.051442819978722
$ws.Cells.Item(16, 9).Value = 1.036801614540814
$ws.Cells.Item(16, 10).Value = 1.032036368024536
$ws.Cells.Item(16, 11).Value = 1.04146202261434
$ws.Cells.Item(16, 12).Value = 1.050162069317522
$ws.Cells.Item(16, 13).Value = 1.054826184920082
$ws.Cells.Item(16, 14).Value = 1.014700098952691

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.025979320877177
$ws.Cells.Item(17, 4).Value = 1.038295586390959
$ws.Cells.Item(17, 5).Value = 1.047098435498968
$ws.Cells.Item(17, 6).Value = 1.051766940670181
$ws.Cells.Item(17, 9).Value = 1.036871485975149
$ws.Cells.Item(17, 10).Value = 1.03222898736589
$ws.Cells.Item(17, 11).Value = 1.041653446018253
$ws.Cells.Item(17, 12).Value = 1.050426048418421
$ws.Cells.Item(17, 13).Value = 1.055078738321232
$ws.Cells.Item(17, 14).Value = 1.014764000559165

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.026172813119937
$ws.Cells.Item(18, 4).Value = 1.03844949461846
$ws.Cells.Item(18, 5).Value = 1.047294521194896
$ws.Cells.Item(18, 6).Value = 1.051956067896676
$ws.Cells.Item(18, 9).Value = 1.036912051062789
$ws.Cells.Item(18, 10).Value = 1.032341312903192
$ws.Cells.Item(18, 11).Value = 1.041765027152689
$ws.Cells.Item(18, 12).Value = 1.050580074530041
$ws.Cells.Item(18, 13).Value = 1.055226044158681
$ws.Cells.Item(18, 14).Value = 1.014801262361008

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.026238797528828
$ws.Cells.Item(19, 4).Value = 1.03850197472334
$ws.Cells.Item(19, 5).Value = 1.047361403034775
$ws.Cells.Item(19, 6).Value = 1.052020567737777
$ws.Cells.Item(19, 9).Value = 1.036925850501645
$ws.Cells.Item(19, 10).Value = 1.032379608552058
$ws.Cells.Item(19, 11).Value = 1.041803061030312
$ws.Cells.Item(19, 12).Value = 1.05063260215409
$ws.Cells.Item(19, 13).Value = 1.055276270867556
$ws.Cells.Item(19, 14).Value = 1.014813965806906

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.025943733504845
$ws.Cells.Item(20, 4).Value = 1.038267276796148
$ws.Cells.Item(20, 5).Value = 1.047062377315679
$ws.Cells.Item(20, 6).Value = 1.05173215801182
$ws.Cells.Item(20, 9).Value = 1.036864009051064
$ws.Cells.Item(20, 10).Value = 1.032208323821173
$ws.Cells.Item(20, 11).Value = 1.04163291562915
$ws.Cells.Item(20, 12).Value = 1.050397720614186
$ws.Cells.Item(20, 13).Value = 1.055051642163726
$ws.Cells.Item(20, 14).Value = 1.014757145648466

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.024985050785953
$ws.Cells.Item(21, 4).Value = 1.037504354211312
$ws.Cells.Item(21, 5).Value = 1.046091726355846
$ws.Cells.Item(21, 6).Value = 1.050795373580428
$ws.Cells.Item(21, 9).Value = 1.036660728321397
$ws.Cells.Item(21, 10).Value = 1.031651156865151
$ws.Cells.Item(21, 11).Value = 1.041078902709081
$ws.Cells.Item(21, 12).Value = 1.049634707553853
$ws.Cells.Item(21, 13).Value = 1.054321305643686
$ws.Cells.Item(21, 14).Value = 1.014572290364819

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.024382755337576
$ws.Cells.Item(22, 4).Value = 1.037024762566438
$ws.Cells.Item(22, 5).Value = 1.045482606839394
$ws.Cells.Item(22, 6).Value = 1.050207051453715
$ws.Cells.Item(22, 9).Value = 1.036531213939039
$ws.Cells.Item(22, 10).Value = 1.031300617064637
$ws.Cells.Item(22, 11).Value = 1.040729925287514
$ws.Cells.Item(22, 12).Value = 1.049155446045356
$ws.Cells.Item(22, 13).Value = 1.053862086197133
$ws.Cells.Item(22, 14).Value = 1.014455968962205

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.024701998864334
$ws.Cells.Item(23, 4).Value = 1.037278994238364
$ws.Cells.Item(23, 5).Value = 1.04580540136803
$ws.Cells.Item(23, 6).Value = 1.050518867882337
$ws.Cells.Item(23, 9).Value = 1.036600033184309
$ws.Cells.Item(23, 10).Value = 1.031486465957705
$ws.Cells.Item(23, 11).Value = 1.040914985751175
$ws.Cells.Item(23, 12).Value = 1.049409466097002
$ws.Cells.Item(23, 13).Value = 1.054105529237551
$ws.Cells.Item(23, 14).Value = 1.014517642056757

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.025959813750086
$ws.Cells.Item(24, 4).Value = 1.03828006865762
$ws.Cells.Item(24, 5).Value = 1.047078670060783
$ws.Cells.Item(24, 6).Value = 1.051747874572244
$ws.Cells.Item(24, 9).Value = 1.036867388137558
$ws.Cells.Item(24, 10).Value = 1.032217660868243
$ws.Cells.Item(24, 11).Value = 1.041642192653655
$ws.Cells.Item(24, 12).Value = 1.050410520570402
$ws.Cells.Item(24, 13).Value = 1.055063885764388
$ws.Cells.Item(24, 14).Value = 1.014760243121344

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.027421421142242
$ws.Cells.Item(25, 4).Value = 1.039442084154535
$ws.Cells.Item(25, 5).Value = 1.048561248861159
$ws.Cells.Item(25, 6).Value = 1.05317692966236
$ws.Cells.Item(25, 9).Value = 1.037170192641684
$ws.Cells.Item(25, 10).Value = 1.033065145114089
$ws.Cells.Item(25, 11).Value = 1.042483211222189
$ws.Cells.Item(25, 12).Value = 1.051574205930246
$ws.Cells.Item(25, 13).Value = 1.056175831244735
$ws.Cells.Item(25, 14).Value = 1.015041338102859

